$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 26-29 (locator values/names refined) ---
$ws.Cells.Item(26, 4).Value = "//li[@id='category-3']/a"

$ws.Cells.Item(27, 2).Value = "AccessoriesMenuLink"
$ws.Cells.Item(27, 4).Value = "//li[@id='category-6']/a"

$ws.Cells.Item(28, 2).Value = "ArtMenuLink"
$ws.Cells.Item(28, 4).Value = "//li[@id='category-9']/a"

$ws.Cells.Item(29, 2).Value = "SubMenuLink"
$ws.Cells.Item(29, 4).Value = "//ul[@class='top-menu']//*[contains(text(),`"{%s}`")]"

# --- Add new "Shopping Method" rows 30-42 ---
$ws.Cells.Item(30, 1).Value = "ShoppingPage"
$ws.Cells.Item(30, 2).Value = "ProductTitleList"
$ws.Cells.Item(30, 3).Value = "xpath"
$ws.Cells.Item(30, 4).Value = "//h2[@class='h3 product-title']/a[text()='{%s}']"

$ws.Cells.Item(31, 1).Value = "ShoppingPage"
$ws.Cells.Item(31, 2).Value = "ProductPriceList"
$ws.Cells.Item(31, 3).Value = "xpath"
$ws.Cells.Item(31, 4).Value = "//h2[@class='h3 product-title']/a//parent::h2//following-sibling::div//span"

$ws.Cells.Item(32, 1).Value = "ShoppingPage"
$ws.Cells.Item(32, 2).Value = "QuanityTxtField"
$ws.Cells.Item(32, 3).Value = "id"
$ws.Cells.Item(32, 4).Value = "quantity_wanted"

$ws.Cells.Item(33, 1).Value = "ShoppingPage"
$ws.Cells.Item(33, 2).Value = "Add2CartBtn"
$ws.Cells.Item(33, 3).Value = "xpath"
$ws.Cells.Item(33, 4).Value = "//button[@class='btn btn-primary add-to-cart']"

$ws.Cells.Item(34, 1).Value = "ShoppingPage"
$ws.Cells.Item(34, 2).Value = "Prcd2ChkOut"
$ws.Cells.Item(34, 3).Value = "xpath"
$ws.Cells.Item(34, 4).Value = "//div[@class='cart-content-btn']/a"

$ws.Cells.Item(35, 1).Value = "ShoppingPage"
$ws.Cells.Item(35, 2).Value = "CartItemCountTxt"
$ws.Cells.Item(35, 3).Value = "xpath"
$ws.Cells.Item(35, 4).Value = "//p[@class='cart-products-count']"

$ws.Cells.Item(36, 1).Value = "ShoppingPage"
$ws.Cells.Item(36, 2).Value = "Prcd2ChkOutConfirmation"
$ws.Cells.Item(36, 3).Value = "xpath"
$ws.Cells.Item(36, 4).Value = "//a[contains(text(),'Proceed to checkout')]"

$ws.Cells.Item(37, 1).Value = "ShoppingPage"
$ws.Cells.Item(37, 2).Value = "AddressConfirmBtn"
$ws.Cells.Item(37, 3).Value = "name"
$ws.Cells.Item(37, 4).Value = "confirm-addresses"

$ws.Cells.Item(38, 1).Value = "ShoppingPage"
$ws.Cells.Item(38, 2).Value = "DlvryConfirmBtn"
$ws.Cells.Item(38, 3).Value = "name"
$ws.Cells.Item(38, 4).Value = "confirmDeliveryOption"

$ws.Cells.Item(39, 1).Value = "ShoppingPage"
$ws.Cells.Item(39, 2).Value = "PaymentOptionCheckbox"
$ws.Cells.Item(39, 3).Value = "xpath"
$ws.Cells.Item(39, 4).Value = "//input[@data-module-name='{%s}']"

$ws.Cells.Item(40, 1).Value = "ShoppingPage"
$ws.Cells.Item(40, 2).Value = "PaymentTextLabel"
$ws.Cells.Item(40, 3).Value = "xpath"
$ws.Cells.Item(40, 4).Value = "//label[contains(@for,'payment-option')]//span"

$ws.Cells.Item(41, 1).Value = "ShoppingPage"
$ws.Cells.Item(41, 2).Value = "PaymentTnCCheckbox"
$ws.Cells.Item(41, 3).Value = "id"
$ws.Cells.Item(41, 4).Value = "conditions_to_approve[terms-and-conditions]"

$ws.Cells.Item(42, 1).Value = "ShoppingPage"
$ws.Cells.Item(42, 2).Value = "PlaceOrderBtn"
$ws.Cells.Item(42, 3).Value = "xpath"
$ws.Cells.Item(42, 4).Value = "//div[@id='payment-confirmation']//button"

# --- Column D width widened to fit the new, longer xpath strings ---
$ws.Columns.Item(4).ColumnWidth = 68.3

# --- View state: scroll/selection mirrors the authored sheet state ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("D48").Select()
